$d = $word.ActiveDocument

# 1) Fix the spelling mistake in the column name: AccountingComapanyID -> AccountingCompanyID
$d.Content.Find.Execute("AccountingComapanyID", $true, $false, $false, $false, $false,
                         $true, 1, $false, "AccountingCompanyID", 2)

# 2) Move the "_GoBack" bookmark from right after the "PackageName ..." paragraph
#    to right after the "Status Bit default 0," paragraph (same CREATE TABLE block).
#    The bookmark must end up collapsed, positioned immediately after that run's text and
#    before the paragraph mark - exactly like it originally was after "PackageName ...".

# Remove the existing (misplaced) bookmark first.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# Find the exact end of the "Status Bit default 0," run.
$found = $d.Content
$found.Find.Execute("Status Bit default 0,", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$endPos = $found.End

# Work around a collapsed-range anchoring quirk in this engine: a zero-length range sitting
# exactly on the character boundary right before a paragraph mark cannot be bookmarked
# directly (it silently anchors to the wrong spot). Temporarily append a marker string
# after the text, bookmark right before the marker, then remove the marker again; the
# bookmark - already anchored as an object - keeps its correct, now-collapsed position.
$marker = "ZZGOBACKMARKERZZ"

$insertPoint = $d.Range($endPos, $endPos)
$insertPoint.InsertAfter($marker)

$withMarker = $d.Content
$withMarker.Find.Execute("Status Bit default 0," + $marker, $true, $false, $false, $false,
                          $false, $true, 1, $false, "", 0)

$bookmarkPos = $withMarker.End - $marker.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$bookmarkRange.Bookmarks.Add("_GoBack")

$markerRange = $d.Range($withMarker.End - $marker.Length, $withMarker.End)
$markerRange.Delete()
